$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row of data appended at the bottom of the log (row 49)
$newRow = 49
$logs.Cells.Item($newRow, 1).Value = "Uitnodiging voor netwerkevent"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Graag nodig ik u uit voor ons zakelijke netwerkevent volgende maand."
$logs.Cells.Item($newRow, 4).Value = "Samenwerking / Partnerverzoek"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 22:36:11"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend conditional formatting ranges to include the new row
$dFc = $logs.Range("D2:D48").FormatConditions
for ($i = 1; $i -le $dFc.Count; $i++) {
    $dFc.Item($i).ModifyAppliesToRange($logs.Range("D2:D49"))
}

$gFc = $logs.Range("G2:G48").FormatConditions
for ($i = 1; $i -le $gFc.Count; $i++) {
    $gFc.Item($i).ModifyAppliesToRange($logs.Range("G2:G49"))
}

# Update the matching Dashboard count (Samenwerking / Partnerverzoek: 13 -> 14)
$dashboard.Cells.Item(2, 2).Value = 14
